# Update countries & provincias Spain
# Applies the 30-Aug-2020 22:11 data refresh to the "Pais" worksheet:
#  - Updates the "last updated" timestamp banner (A1)
#  - Refreshes total/new/active/recovered/critical/death stats for several
#    countries whose case counts changed
#  - Ruanda's new-case surge pushes it above Congo/Cuba/Surinam/Eslovaquia
#    in the ranking, so those four rows shift down one position
#  - Togo's new-case surge pushes it above Letonia, so that row swaps too

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Agosto de 2020 a las 22:11"

# --- Helper: write B:H (Casos totales..Muertes) for a row -------------
function Set-Stats($row, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Range("B$row").Value = $total
    $ws.Range("C$row").Value = $nuevos
    $ws.Range("D$row").Value = $activos
    $ws.Range("E$row").Value = $recuperados
    $ws.Range("F$row").Value = $criticos
    $ws.Range("G$row").Value = $muertesHoy
    $ws.Range("H$row").Value = $muertes
}

# Row 4 - Estados Unidos
Set-Stats 4 6163159 23789 3415980 2560079 0 245 187100

# Row 23 - Alemania
Set-Stats 23 243282 457 217484 16434 0 1 9364

# Row 79 - Costa de Marfil
Set-Stats 79 17948 55 16553 1280 0 0 115

# Row 96 - Guinea
Set-Stats 96 9309 58 8333 917 0 0 59

# Row 104 - Mauritania
Set-Stats 104 7016 4 6430 427 0 1 159

# Row 106 - Zimbabue
Set-Stats 106 6412 6 5061 1155 0 0 196

# Row 108 - Malaui
Set-Stats 108 5536 8 3147 2215 0 0 174

# Row 114 - Suazilandia (stats update only, stays in place)
Set-Stats 114 4561 51 3478 992 0 0 91

# Rows 116-120: Ruanda jumps ahead of Congo/Cuba/Surinam/Eslovaquia.
# Row 116 becomes Ruanda with its fresh stats; Congo, Cuba, Surinam and
# Eslovaquia each shift down one row, keeping their own (unchanged) stats.
$ws.Range("A116").Value = "Ruanda"
Set-Stats 116 4020 177 1918 2086 0 0 16

$ws.Range("A117").Value = "Congo"
Set-Stats 117 3979 0 1742 2159 0 0 78

$ws.Range("A118").Value = "Cuba"
Set-Stats 118 3973 48 3327 552 0 0 94

$ws.Range("A119").Value = "Surinam"
Set-Stats 119 3954 0 2991 896 0 0 67

$ws.Range("A120").Value = "Eslovaquia"
Set-Stats 120 3876 34 2278 1565 0 0 33

# Row 125 - Somalia
Set-Stats 125 3310 0 2481 731 0 1 98

# Row 132 - Mali
Set-Stats 132 2773 16 2169 478 0 0 126

# Row 134 - Angola
Set-Stats 134 2624 73 1063 1454 0 0 107

# Row 146 - Malta
Set-Stats 146 1862 15 1313 537 0 1 12

# Rows 154-155: Togo jumps ahead of Letonia.
# Row 154 becomes Togo with its fresh stats; Letonia shifts down to row
# 155, keeping its own (unchanged) stats.
$ws.Range("A154").Value = "Togo"
Set-Stats 154 1396 6 996 373 0 0 27

$ws.Range("A155").Value = "Letonia"
Set-Stats 155 1393 12 1163 196 0 0 34

# Row 166 - Santo Tome y Principe
Set-Stats 166 896 1 851 30 0 0 15
